$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '65.080.94'
Set-TextValue 'E2' '  +2.04%  '
Set-TextValue 'D3' '3.190.27'
Set-TextValue 'E3' '  +1.50%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '595.26'
Set-TextValue 'E5' '  +1.43%  '
Set-TextValue 'D6' '153.95'
Set-TextValue 'E6' '  +5.75%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '3.188.49'
Set-TextValue 'E8' '  +1.59%  '
Set-TextValue 'D9' '0.541'
Set-TextValue 'E9' '  +2.30%  '
Set-TextValue 'D10' '0.165'
Set-TextValue 'E10' '  +2.65%  '
Set-TextValue 'D11' '6.03'
Set-TextValue 'E11' '  +5.27%  '
Set-TextValue 'D12' '0.470'
Set-TextValue 'E12' '  +2.81%  '
Set-TextValue 'D13' '0.0000254'
Set-TextValue 'E13' '  +2.82%  '
Set-TextValue 'D14' '39.29'
Set-TextValue 'E14' '  +6.16%  '
Set-TextValue 'D15' '3.718.99'
Set-TextValue 'E15' '  +1.51%  '
Set-TextValue 'E16' '  +0.21%  '
Set-TextValue 'D17' '7.41'
Set-TextValue 'E17' '  +4.81%  '
Set-TextValue 'D18' '64.843.67'
Set-TextValue 'E18' '  +1.94%  '
Set-TextValue 'D19' '3.194.52'
Set-TextValue 'E19' '  +1.59%  '
Set-TextValue 'D20' '481.41'
Set-TextValue 'E20' '  +3.76%  '
Set-TextValue 'D21' '15.07'
Set-TextValue 'E21' '  +5.65%  '
Set-TextValue 'D22' '0.769'
Set-TextValue 'E22' '  +5.26%  '
Set-TextValue 'D23' '7.86'
Set-TextValue 'E23' '  +5.75%  '
Set-TextValue 'D24' '13.72'
Set-TextValue 'E24' '  +5.84%  '
Set-TextValue 'D25' '2.45'
Set-TextValue 'E25' '  +10.62%  '
Set-TextValue 'D26' '83.59'
Set-TextValue 'E26' '  +2.89%  '
Set-TextValue 'B27' 'Dai'
Set-TextValue 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D27' '1.00'
Set-TextValue 'E27' '  +0.16%  '
Set-TextValue 'B28' 'RenderToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D28' '10.00'
Set-TextValue 'E28' '  +8.26%  '
Set-TextValue 'D29' '2.77'
Set-TextValue 'E29' '  +3.21%  '
Set-TextValue 'D30' '7.49'
Set-TextValue 'E30' '  +7.60%  '
Set-TextValue 'D31' '2.27'
Set-TextValue 'E31' '  +2.52%  '
Set-TextValue 'E32' '  +0.33%  '
Set-TextValue 'D33' '0.119'
Set-TextValue 'E33' '  +8.50%  '
Set-TextValue 'D34' '28.40'
Set-TextValue 'E34' '  +5.22%  '
Set-TextValue 'D35' '0.0₃0903'
Set-TextValue 'E35' '  +6.35%  '
Set-TextValue 'D36' '3.60'
Set-TextValue 'E36' '  +8.53%  '
Set-TextValue 'D37' '1.09'
Set-TextValue 'E37' '  +4.48%  '
Set-TextValue 'B38' 'Filecoin'
Set-TextValue 'C38' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D38' '6.33'
Set-TextValue 'E38' '  +5.35%  '
Set-TextValue 'B39' 'Stacks'
Set-TextValue 'C39' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '2.37'
Set-TextValue 'E39' '  +3.41%  '
Set-TextValue 'D40' '475.28'
Set-TextValue 'E40' '  +7.82%  '
Set-TextValue 'D41' '9.44'
Set-TextValue 'E41' '  +7.04%  '
Set-TextValue 'D42' '51.67'
Set-TextValue 'E42' '  +1.37%  '
Set-TextValue 'D43' '0.303'
Set-TextValue 'E43' '  +8.72%  '
Set-TextValue 'D44' '0.0382'
Set-TextValue 'E44' '  +2.87%  '
Set-TextValue 'D45' '2.950.43'
Set-TextValue 'E45' '  +1.68%  '
Set-TextValue 'E46' '  +4.26%  '
Set-TextValue 'D47' '38.91'
Set-TextValue 'E47' '  +4.88%  '
Set-TextValue 'D48' '131.62'
Set-TextValue 'E48' '  +4.75%  '
Set-TextValue 'B49' 'InjectiveProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D49' '25.97'
Set-TextValue 'E49' '  +6.91%  '
Set-TextValue 'B50' 'ThetaToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D50' '2.34'
Set-TextValue 'E50' '  +7.96%  '
Set-TextValue 'E51' '  +0.01%  '
